$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new StatQuery text (Cypher query) with embedded newlines.
$nl = [char]10
$newQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)" + $nl +
    "OPTIONAL MATCH (samp:sample)-->(c)" + $nl +
    "OPTIONAL MATCH (diag:diagnosis)-->(c)" + $nl +
    "OPTIONAL MATCH (f:file)-[*]->(c)" + $nl +
    "OPTIONAL MATCH (sf:file)-->(s)" + $nl +
    "WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p" + $nl +
    "WHERE demo.breed IN ['German Shorthaired Pointer']" + $nl +
    "RETURN  " + $nl +
    "    count(distinct p) AS Programs," + $nl +
    "    count(distinct s) AS Studies," + $nl +
    "    count(distinct c) AS Cases," + $nl +
    "    count(distinct samp) AS Samples," + $nl +
    "    count(distinct f) AS ``Case Files``," + $nl +
    "    count(distinct sf) AS ``Study Files``"

# Update the StatQuery column (C) for all three rows so the shared string
# is rewritten in place instead of orphaned.
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Adjust row heights (previously capped at 409.6, now shorter to fit new text).
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Update view: zoom level and selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$ws.Range("B4:B5").Select()
